# Update the "Primary Review Worksheet" header row (row 1) to the new
# RAI-oriented column layout: a few existing headers shift position, the
# unused funding-by-year columns (Total/BY1-5 Funding, Has Keywords) are
# dropped, and a block of new "RAI ..." / "POC ..." review columns is
# appended, extending the sheet from column AA out to column AG.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Primary Review Worksheet")

$ws.Range("A1").Value  = "Primary Reviewer"
$ws.Range("B1").Value  = "AI Analysis"
$ws.Range("C1").Value  = "FY22 Label"
$ws.Range("D1").Value  = "Service/DoD Component Reviewer"
$ws.Range("E1").Value  = "FY22 Service Reviewer"
$ws.Range("F1").Value  = "Planned Transition Partner"
$ws.Range("G1").Value  = "Current Mission Partners (Academia, Industry, or Other)"
$ws.Range("H1").Value  = "Primary Reviewer Notes"
$ws.Range("I1").Value  = "FY22 POC Reviewer"
$ws.Range("J1").Value  = "FY (BY1)"
$ws.Range("K1").Value  = "Doc Type"
$ws.Range("L1").Value  = "Service / Agency"
$ws.Range("M1").Value  = "Agency / Office"
$ws.Range("N1").Value  = "APPN Symbol"
$ws.Range("O1").Value  = "APPN Title"
$ws.Range("P1").Value  = "Project"
$ws.Range("Q1").Value  = "BA"
$ws.Range("R1").Value  = "BA Title"
$ws.Range("S1").Value  = "PE / BLI"
$ws.Range("T1").Value  = "Project # (RDT&E Only)"
$ws.Range("U1").Value  = "Department"
$ws.Range("V1").Value  = "RAI Secondary Reviewer"
$ws.Range("W1").Value  = "RAI Tag Agree"
$ws.Range("X1").Value  = "RAI Tag"
$ws.Range("Y1").Value  = "RAI Transition Partner Agree"
$ws.Range("Z1").Value  = "RAI Transition Partner"
$ws.Range("AA1").Value = "RAI Mission Partners"
$ws.Range("AB1").Value = "POC Title"
$ws.Range("AC1").Value = "POC Name"
$ws.Range("AD1").Value = "POC Email"
$ws.Range("AE1").Value = "POC Org"
$ws.Range("AF1").Value = "POC Phone Number"
$ws.Range("AG1").Value = "RAI Review Notes"

# Restore the author's on-disk selection/viewport (cell X18 selected).
$ws.Range("X18").Select() | Out-Null
